$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# swap F:V between row 257 and row 258
$ws.Range("F257").Value2 = "ABC"
$ws.Range("F258").Value2 = "Ceara"
$ws.Range("G257").Value2 = 1
$ws.Range("G258").Value2 = 1
$ws.Range("H257").Value2 = "Sampaio Correa"
$ws.Range("H258").Value2 = "Criciuma"
$ws.Range("I257").Value2 = 1
$ws.Range("I258").Value2 = 0
$ws.Range("J257").Value2 = 2.56
$ws.Range("J258").Value2 = 2.01
$ws.Range("K257").Value2 = "29/08/2023 00:12"
$ws.Range("K258").Value2 = "27/08/2023 16:13"
$ws.Range("L257").Value2 = 2.52
$ws.Range("L258").Value2 = 1.98
$ws.Range("M257").Value2 = "02/09/2023 21:59"
$ws.Range("M258").Value2 = "02/09/2023 21:57"
$ws.Range("N257").Value2 = 2.84
$ws.Range("N258").Value2 = 3.17
$ws.Range("O257").Value2 = "29/08/2023 00:12"
$ws.Range("O258").Value2 = "27/08/2023 16:13"
$ws.Range("P257").Value2 = 2.78
$ws.Range("P258").Value2 = 3.28
$ws.Range("Q257").Value2 = "02/09/2023 21:59"
$ws.Range("Q258").Value2 = "02/09/2023 21:55"
$ws.Range("R257").Value2 = 3.16
$ws.Range("R258").Value2 = 4.02
$ws.Range("S257").Value2 = "29/08/2023 00:12"
$ws.Range("S258").Value2 = "27/08/2023 16:13"
$ws.Range("T257").Value2 = 3.61
$ws.Range("T258").Value2 = 4.45
$ws.Range("U257").Value2 = "02/09/2023 21:51"
$ws.Range("U258").Value2 = "02/09/2023 21:57"
$ws.Range("V257").Value2 = "https://www.betexplorer.com/football/brazil/serie-b/abc-sampaio-correa/dpmoRnrq/"
$ws.Range("V258").Value2 = "https://www.betexplorer.com/football/brazil/serie-b/ceara-criciuma/p67AXl5F/"

# swap F:V between row 276 and row 277
$ws.Range("F276").Value2 = "Botafogo SP"
$ws.Range("F277").Value2 = "Sampaio Correa"
$ws.Range("G276").Value2 = 1
$ws.Range("G277").Value2 = 2
$ws.Range("H276").Value2 = "Atletico GO"
$ws.Range("H277").Value2 = "Chapecoense-SC"
$ws.Range("I276").Value2 = 0
$ws.Range("I277").Value2 = 0
$ws.Range("J276").Value2 = 2.92
$ws.Range("J277").Value2 = 2.08
$ws.Range("K276").Value2 = "11/09/2023 08:12"
$ws.Range("K277").Value2 = "11/09/2023 02:42"
$ws.Range("L276").Value2 = 3.52
$ws.Range("L277").Value2 = 2.12
$ws.Range("M276").Value2 = "16/09/2023 21:52"
$ws.Range("M277").Value2 = "16/09/2023 21:50"
$ws.Range("N276").Value2 = 2.79
$ws.Range("N277").Value2 = 2.96
$ws.Range("O276").Value2 = "11/09/2023 08:12"
$ws.Range("O277").Value2 = "11/09/2023 02:42"
$ws.Range("P276").Value2 = 2.88
$ws.Range("P277").Value2 = 2.89
$ws.Range("Q276").Value2 = "16/09/2023 21:50"
$ws.Range("Q277").Value2 = "16/09/2023 21:50"
$ws.Range("R276").Value2 = 2.92
$ws.Range("R277").Value2 = 4.43
$ws.Range("S276").Value2 = "11/09/2023 08:12"
$ws.Range("S277").Value2 = "11/09/2023 02:42"
$ws.Range("T276").Value2 = 2.48
$ws.Range("T277").Value2 = 4.62
$ws.Range("U276").Value2 = "16/09/2023 21:50"
$ws.Range("U277").Value2 = "16/09/2023 21:50"
$ws.Range("V276").Value2 = "https://www.betexplorer.com/football/brazil/serie-b/botafogo-sp-atletico-go/M5cNa7l8/"
$ws.Range("V277").Value2 = "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-chapecoense-sc/6gaZd5YQ/"

# swap F:V between row 286 and row 287
$ws.Range("F286").Value2 = "Chapecoense-SC"
$ws.Range("F287").Value2 = "CRB"
$ws.Range("G286").Value2 = 1
$ws.Range("G287").Value2 = 1
$ws.Range("H286").Value2 = "Ceara"
$ws.Range("H287").Value2 = "Guarani"
$ws.Range("I286").Value2 = 1
$ws.Range("I287").Value2 = 0
$ws.Range("J286").Value2 = 2.55
$ws.Range("J287").Value2 = 2.09
$ws.Range("K286").Value2 = "19/09/2023 01:12"
$ws.Range("K287").Value2 = "18/09/2023 20:13"
$ws.Range("L286").Value2 = 2.4
$ws.Range("L287").Value2 = 2.14
$ws.Range("M286").Value2 = "23/09/2023 21:58"
$ws.Range("M287").Value2 = "23/09/2023 21:53"
$ws.Range("N286").Value2 = 2.89
$ws.Range("N287").Value2 = 3.03
$ws.Range("O286").Value2 = "19/09/2023 01:12"
$ws.Range("O287").Value2 = "18/09/2023 20:13"
$ws.Range("P286").Value2 = 2.91
$ws.Range("P287").Value2 = 3.06
$ws.Range("Q286").Value2 = "23/09/2023 21:59"
$ws.Range("Q287").Value2 = "23/09/2023 21:53"
$ws.Range("R286").Value2 = 3.27
$ws.Range("R287").Value2 = 4.27
$ws.Range("S286").Value2 = "19/09/2023 01:12"
$ws.Range("S287").Value2 = "18/09/2023 20:13"
$ws.Range("T286").Value2 = 3.67
$ws.Range("T287").Value2 = 4.15
$ws.Range("U286").Value2 = "23/09/2023 21:59"
$ws.Range("U287").Value2 = "23/09/2023 21:53"
$ws.Range("V286").Value2 = "https://www.betexplorer.com/football/brazil/serie-b/chapecoense-sc-ceara/OfNy52Qm/"
$ws.Range("V287").Value2 = "https://www.betexplorer.com/football/brazil/serie-b/crb-guarani/KMGp3te0/"

# swap F:V between row 304 and row 305
$ws.Range("F304").Value2 = "CRB"
$ws.Range("F305").Value2 = "Atletico GO"
$ws.Range("G304").Value2 = 2
$ws.Range("G305").Value2 = 1
$ws.Range("H304").Value2 = "Ceara"
$ws.Range("H305").Value2 = "Ituano"
$ws.Range("I304").Value2 = 0
$ws.Range("I305").Value2 = 0
$ws.Range("J304").Value2 = 2.2
$ws.Range("J305").Value2 = 1.66
$ws.Range("K304").Value2 = "01/10/2023 22:12"
$ws.Range("K305").Value2 = "01/10/2023 22:12"
$ws.Range("L304").Value2 = 2.01
$ws.Range("L305").Value2 = 1.63
$ws.Range("M304").Value2 = "07/10/2023 02:24"
$ws.Range("M305").Value2 = "07/10/2023 02:26"
$ws.Range("N304").Value2 = 3.02
$ws.Range("N305").Value2 = 3.63
$ws.Range("O304").Value2 = "01/10/2023 22:12"
$ws.Range("O305").Value2 = "01/10/2023 22:12"
$ws.Range("P304").Value2 = 3.26
$ws.Range("P305").Value2 = 3.77
$ws.Range("Q304").Value2 = "07/10/2023 02:24"
$ws.Range("Q305").Value2 = "07/10/2023 02:26"
$ws.Range("R304").Value2 = 3.87
$ws.Range("R305").Value2 = 5.97
$ws.Range("S304").Value2 = "01/10/2023 22:12"
$ws.Range("S305").Value2 = "01/10/2023 22:12"
$ws.Range("T304").Value2 = 4.36
$ws.Range("T305").Value2 = 6.24
$ws.Range("U304").Value2 = "07/10/2023 02:24"
$ws.Range("U305").Value2 = "07/10/2023 02:26"
$ws.Range("V304").Value2 = "https://www.betexplorer.com/football/brazil/serie-b/crb-ceara/IgpBFApN/"
$ws.Range("V305").Value2 = "https://www.betexplorer.com/football/brazil/serie-b/atletico-go-ituano/OWOy1BVT/"

# swap F:V between row 344 and row 345
$ws.Range("F344").Value2 = "Juventude"
$ws.Range("F345").Value2 = "Londrina"
$ws.Range("G344").Value2 = 2
$ws.Range("G345").Value2 = 0
$ws.Range("H344").Value2 = "Ituano"
$ws.Range("H345").Value2 = "Guarani"
$ws.Range("I344").Value2 = 1
$ws.Range("I345").Value2 = 0
$ws.Range("J344").Value2 = 1.77
$ws.Range("J345").Value2 = 2.84
$ws.Range("K344").Value2 = "29/10/2023 23:12"
$ws.Range("K345").Value2 = "29/10/2023 22:42"
$ws.Range("L344").Value2 = 1.83
$ws.Range("L345").Value2 = 2.88
$ws.Range("M344").Value2 = "04/11/2023 01:12"
$ws.Range("M345").Value2 = "04/11/2023 01:25"
$ws.Range("N344").Value2 = 3.34
$ws.Range("N345").Value2 = 2.72
$ws.Range("O344").Value2 = "29/10/2023 23:12"
$ws.Range("O345").Value2 = "29/10/2023 22:42"
$ws.Range("P344").Value2 = 3.33
$ws.Range("P345").Value2 = 2.99
$ws.Range("Q344").Value2 = "04/11/2023 01:12"
$ws.Range("Q345").Value2 = "04/11/2023 01:25"
$ws.Range("R344").Value2 = 5.52
$ws.Range("R345").Value2 = 2.99
$ws.Range("S344").Value2 = "29/10/2023 23:12"
$ws.Range("S345").Value2 = "29/10/2023 22:42"
$ws.Range("T344").Value2 = 5.27
$ws.Range("T345").Value2 = 2.85
$ws.Range("U344").Value2 = "04/11/2023 01:22"
$ws.Range("U345").Value2 = "04/11/2023 01:25"
$ws.Range("V344").Value2 = "https://www.betexplorer.com/football/brazil/serie-b/esporte-clube-juventude-ituano/Ug78kHRg/"
$ws.Range("V345").Value2 = "https://www.betexplorer.com/football/brazil/serie-b/londrina-guarani/SvQBlys0/"

# swap F:V between row 346 and row 347
$ws.Range("F346").Value2 = "Ponte Preta"
$ws.Range("F347").Value2 = "Botafogo SP"
$ws.Range("G346").Value2 = 0
$ws.Range("G347").Value2 = 2
$ws.Range("H346").Value2 = "Avai"
$ws.Range("H347").Value2 = "Ceara"
$ws.Range("I346").Value2 = 1
$ws.Range("I347").Value2 = 2
$ws.Range("J346").Value2 = 2.33
$ws.Range("J347").Value2 = 2.41
$ws.Range("K346").Value2 = "29/10/2023 02:42"
$ws.Range("K347").Value2 = "29/10/2023 04:42"
$ws.Range("L346").Value2 = 2.52
$ws.Range("L347").Value2 = 2.28
$ws.Range("M346").Value2 = "04/11/2023 20:45"
$ws.Range("M347").Value2 = "04/11/2023 20:59"
$ws.Range("N346").Value2 = 2.97
$ws.Range("N347").Value2 = 2.88
$ws.Range("O346").Value2 = "29/10/2023 02:42"
$ws.Range("O347").Value2 = "29/10/2023 04:42"
$ws.Range("P346").Value2 = 2.86
$ws.Range("P347").Value2 = 2.98
$ws.Range("Q346").Value2 = "04/11/2023 20:38"
$ws.Range("Q347").Value2 = "04/11/2023 20:59"
$ws.Range("R346").Value2 = 3.61
$ws.Range("R347").Value2 = 3.55
$ws.Range("S346").Value2 = "29/10/2023 02:42"
$ws.Range("S347").Value2 = "29/10/2023 04:42"
$ws.Range("T346").Value2 = 3.48
$ws.Range("T347").Value2 = 3.86
$ws.Range("U346").Value2 = "04/11/2023 20:45"
$ws.Range("U347").Value2 = "04/11/2023 20:59"
$ws.Range("V346").Value2 = "https://www.betexplorer.com/football/brazil/serie-b/ponte-preta-avai/SQpwzxS6/"
$ws.Range("V347").Value2 = "https://www.betexplorer.com/football/brazil/serie-b/botafogo-sp-ceara/xdSTdaZJ/"

# Add new row 359 (copy formatting from row 358, then set values)
$ws.Range("A358:V358").Copy()
$ws.Range("A359").PasteSpecial(-4122)

$ws.Range("A359").Value2 = 358
$ws.Range("B359").Value2 = "brazil"
$ws.Range("C359").Value2 = "serie-b"
$ws.Range("D359").NumberFormat = "@"
$ws.Range("D359").Value2 = "2023"
$ws.Range("D359").Style = "Normal"
$ws.Range("E359").Value2 = 45242.91666666666
$ws.Range("F359").Value2 = "Novorizontino"
$ws.Range("G359").Value2 = 1
$ws.Range("H359").Value2 = "Vitoria"
$ws.Range("I359").Value2 = 2
$ws.Range("J359").Value2 = 1.73
$ws.Range("K359").Value2 = "05/11/2023 22:12"
$ws.Range("L359").Value2 = 1.91
$ws.Range("M359").Value2 = "12/11/2023 21:55"
$ws.Range("N359").Value2 = 3.42
$ws.Range("O359").Value2 = "05/11/2023 22:12"
$ws.Range("P359").Value2 = 3.23
$ws.Range("Q359").Value2 = "12/11/2023 21:55"
$ws.Range("R359").Value2 = 5.71
$ws.Range("S359").Value2 = "05/11/2023 22:12"
$ws.Range("T359").Value2 = 4.96
$ws.Range("U359").Value2 = "12/11/2023 21:58"
$ws.Range("V359").Value2 = "https://www.betexplorer.com/football/brazil/serie-b/novorizontino-vitoria/xEZzk5ht/"
